# Updates crypto price/volume figures per the scraper run on
# Thu Jan 19 08:40:52 UTC 2023 (GitHub Actions symbol-list refresh).
# D = Price, E = Volume(1h) change; both columns are plain text cells
# (values such as "292.03" or "-3.15%" are literal strings, not
# numbers/percentages), so NumberFormat is forced to "@" (Text)
# before assignment to stop Excel from auto-coercing the numeric-
# looking / percent-looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "292.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.15%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.56%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.962"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.44%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07233"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.36%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.791"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-7.09%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.682"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.95%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.762"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.20%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8976"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.50%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1664"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.51%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07712"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.39%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07977"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.44%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03040"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.94%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.17%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001501"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.08%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005752"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.85%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.466"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.19%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.084"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.26%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.87%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.88%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.054"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.27%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "13.06%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04508"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.90%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001217"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.71%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004014"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-9.06%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001252"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.01%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01598"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.36%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04415"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.60%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007242"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-5.65%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1308"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.08%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007683"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-18.16%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009512"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.67%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005929"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.38%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "173.66%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003004"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.25%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
